# Re-doing global M2 module: refresh M2/FX length counters and
# first/last-date watermarks in the Global M2 DataComp sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - China
$ws.Range("C2").Value = 360
$ws.Range("F2").Value = "12/01/2025"
$ws.Range("G2").Value = "07/02/1984"
$ws.Range("H2").Value = "02/02/2026"

# Row 3 - United States
$ws.Range("E3").Value = "05/01/1984"
$ws.Range("F3").Value = "12/01/2025"

# Row 4 - Euro Area
$ws.Range("E4").Value = "05/01/1984"
$ws.Range("F4").Value = "12/01/2025"
$ws.Range("G4").Value = "07/02/1984"
$ws.Range("H4").Value = "02/02/2026"

# Row 5 - Japan
$ws.Range("E5").Value = "05/01/1984"
$ws.Range("F5").Value = "12/01/2025"
$ws.Range("G5").Value = "07/02/1984"
$ws.Range("H5").Value = "02/02/2026"

# Row 6 - United Kingdom
$ws.Range("G6").Value = "07/02/1984"
$ws.Range("H6").Value = "02/02/2026"

# Row 7 - South Korea
$ws.Range("E7").Value = "04/01/1984"
$ws.Range("F7").Value = "11/01/2025"
$ws.Range("G7").Value = "07/02/1984"
$ws.Range("H7").Value = "02/02/2026"

# Row 8 - Hong Kong
$ws.Range("D8").Value = 436
$ws.Range("H8").Value = "02/02/2026"

# Row 9 - Australia
$ws.Range("G9").Value = "07/02/1984"
$ws.Range("H9").Value = "02/02/2026"

# Row 10 - Taiwan
$ws.Range("E10").Value = "05/01/1984"
$ws.Range("F10").Value = "12/01/2025"
$ws.Range("G10").Value = "07/02/1984"
$ws.Range("H10").Value = "02/02/2026"

# Row 11 - Canada
$ws.Range("E11").Value = "04/01/1984"
$ws.Range("F11").Value = "11/01/2025"
$ws.Range("G11").Value = "07/02/1984"
$ws.Range("H11").Value = "02/02/2026"

# Row 12 - Russia
$ws.Range("C12").Value = 397
$ws.Range("D12").Value = 378
$ws.Range("F12").Value = "12/01/2025"
$ws.Range("H12").Value = "02/02/2026"

# Row 13 - Switzerland
$ws.Range("C13").Value = 493
$ws.Range("F13").Value = "12/01/2025"
$ws.Range("G13").Value = "07/02/1984"
$ws.Range("H13").Value = "02/02/2026"

# Row 14 - Brazil
$ws.Range("C14").Value = 448
$ws.Range("D14").Value = 422
$ws.Range("F14").Value = "11/01/2025"
$ws.Range("H14").Value = "02/02/2026"

# Row 15 - India
$ws.Range("C15").Value = 409
$ws.Range("F15").Value = "11/01/2025"
$ws.Range("G15").Value = "06/01/1984"
$ws.Range("H15").Value = "02/02/2026"

# Row 16 - Mexico
$ws.Range("D16").Value = 436
$ws.Range("H16").Value = "02/02/2026"

# Row 17 - Saudi Arabia
$ws.Range("C17").Value = 396
$ws.Range("D17").Value = 420
$ws.Range("F17").Value = "12/01/2025"
$ws.Range("H17").Value = "02/02/2026"

# Row 18 - Singapore
$ws.Range("D18").Value = 286
$ws.Range("E18").Value = "05/01/1984"
$ws.Range("F18").Value = "12/01/2025"
$ws.Range("H18").Value = "02/02/2026"

# Row 19 - Indonesia
$ws.Range("D19").Value = 424
$ws.Range("E19").Value = "05/01/1984"
$ws.Range("F19").Value = "12/01/2025"
$ws.Range("H19").Value = "02/02/2026"

# Row 20 - Malaysia
$ws.Range("E20").Value = "12/01/1977"
$ws.Range("F20").Value = "12/01/2025"
$ws.Range("G20").Value = "07/02/1984"
$ws.Range("H20").Value = "02/02/2026"

# Row 21 - Norway
$ws.Range("E21").Value = "05/01/1984"
$ws.Range("F21").Value = "12/01/2025"
$ws.Range("G21").Value = "07/02/1984"
$ws.Range("H21").Value = "02/02/2026"

# Row 22 - Philippines
$ws.Range("D22").Value = 406
$ws.Range("E22").Value = "04/01/1984"
$ws.Range("F22").Value = "11/01/2025"
$ws.Range("H22").Value = "02/02/2026"

# Row 23 - New Zealand
$ws.Range("D23").Value = 243
$ws.Range("H23").Value = "02/02/2026"

# Row 24 - Denmark
$ws.Range("C24").Value = 420
$ws.Range("F24").Value = "12/01/2025"
$ws.Range("G24").Value = "07/02/1984"
$ws.Range("H24").Value = "02/02/2026"

# Row 25 - South Africa
$ws.Range("E25").Value = "05/01/1984"
$ws.Range("F25").Value = "12/01/2025"
$ws.Range("G25").Value = "07/02/1984"
$ws.Range("H25").Value = "02/02/2026"

# Row 27 - Colombia
$ws.Range("D27").Value = 436
$ws.Range("E27").Value = "05/01/1984"
$ws.Range("F27").Value = "12/01/2025"
$ws.Range("H27").Value = "02/02/2026"

# Row 28 - Kuwait
$ws.Range("C28").Value = 385
$ws.Range("D28").Value = 393
$ws.Range("F28").Value = "12/01/2025"
$ws.Range("H28").Value = "02/02/2026"
